# CASMNET-858: rename the INTER_SWITCH_LINKS tab to SWITCH_TO_SWITCH and
# update the workbook's "last looked at" UI state to match: the renamed
# sheet becomes the active/selected tab (with a new cell selection),
# while the sheet that used to be active (COMPUTE_NODES) is no longer
# the selected tab.

$wb = $excel.ActiveWorkbook

# Rename the sheet.
$ws = $wb.Worksheets.Item("INTER_SWITCH_LINKS")
$ws.Name = "SWITCH_TO_SWITCH"

# Make it the active sheet/tab, and move the cell selection on it.
$ws.Activate()
$ws.Range("E29").Select()
